$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.123.72"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.790.22"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'229.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'32.53"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "'0.289"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "'0.0937"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "2.048.62"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'11.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "1.800.87"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "34.077.64"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "'68.45"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "'245.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "'160.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").Value = "'16.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("D31").Value = "'0.0515"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").Value = "'3.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "1.398.97"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").Value = "'0.661"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.54%  "
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").Value = "'0.0188"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "'2.22"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'78.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "'13.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.11%  "
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "'109.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("E47").Value = "  +9.26%  "
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Value = "1.947.75"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  +0.27%  "
